$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colors (OLE BGR-packed values)
$lightBlue = 15773696   # FF00B0F0 - new "done" highlight color
$yellow    = 65535      # FFFFFF00 - existing "pending" color
$green     = 5287936    # FF00B050 - existing "reviewed, no change" color
$darkBlue  = 6299648    # FF002060 - existing "N/A" color

# ---------------------------------------------------------------------
# Rows 2-16: upgrade the existing review-status highlighting in columns
# D and E from the old yellow/green scheme to the new light-blue scheme.
# (Text/values in these cells are unchanged - only fill colors change.)
# ---------------------------------------------------------------------
$ws.Range("D2:D4").Interior.Color = $lightBlue
$ws.Range("D9:D12").Interior.Color = $lightBlue
$ws.Range("D13:D16").Interior.Color = $lightBlue

$ws.Range("E2").Interior.Color = $lightBlue
$ws.Range("E3:E4").Interior.Color = $lightBlue
$ws.Range("E9:E12").Interior.Color = $lightBlue
$ws.Range("E13:E15").Interior.Color = $lightBlue

# ---------------------------------------------------------------------
# Rows 19-25: these were previously blank in columns D/E; fill them in
# with review results, using the pre-existing yellow/green/dark-blue
# highlight scheme (the same one rows 2-16 used before being upgraded).
# Row 17/24 also pick up new notes in columns F/G.
# (Order below matches the order new shared strings were authored in.)
# ---------------------------------------------------------------------
$ws.Range("D19").Value = "4 revisions: 3 insertions, 1 deletion"
$ws.Range("D19").Interior.Color = $yellow
$ws.Range("E19").Value = "review complete - no change needed"
$ws.Range("E19").Interior.Color = $green

$ws.Range("D20").Value = "4 revisions: 2 insertions, 2 deletions"
$ws.Range("D20").Interior.Color = $yellow
$ws.Range("E20").Value = "review complete - no change needed"
$ws.Range("E20").Interior.Color = $green

$ws.Range("D21").Value = "27 revisions: 21 insertions, 6 deletions"
$ws.Range("D21").Interior.Color = $yellow
$ws.Range("E21").Value = "4 revisions: 4 insertions, 0 deletions"
$ws.Range("E21").Interior.Color = $yellow

$ws.Range("D22").Value = "31 revisions: 23 insertions, 8 deletions"
$ws.Range("D22").Interior.Color = $yellow
$ws.Range("E22").Value = "7 revisions: 5 insertions, 1 deletions"
$ws.Range("E22").Interior.Color = $yellow

# Row 17: new note in column F
$ws.Range("F17").Value = "Material_property_definition_schema"

$ws.Range("D23").Value = "4 revisions: 3 insertions, 1 deletion"
$ws.Range("D23").Interior.Color = $yellow
$ws.Range("E23").Value = "review complete - no change needed"
$ws.Range("E23").Interior.Color = $green

# Row 24: also gets two new notes in columns F and G
$ws.Range("F24").Value = "State_observed_schema"
$ws.Range("G24").Value = "State_schema"
$ws.Range("D24").Value = "6 revisions: 3 insertions, 3 deletions"
$ws.Range("D24").Interior.Color = $yellow
$ws.Range("E24").Value = "19 revisions: 11 insertions, 8 deletions"
$ws.Range("E24").Interior.Color = $yellow

$ws.Range("D25").Value = "10 revisions: 8 insertions, 2 deletions"
$ws.Range("D25").Interior.Color = $yellow
# E25 reuses the existing dark-blue/white-font "N/A" style (fontId 5 is a
# theme-white color, so copy the format from E17 rather than re-setting
# Font.Color, which would otherwise create a near-duplicate font entry).
$ws.Range("E17").Copy()
$ws.Range("E25").PasteSpecial(-4122)
$ws.Range("E25").Value = "N/A"

# ---------------------------------------------------------------------
# Selection moves to D25 (matches the author's final cursor position)
# ---------------------------------------------------------------------
$ws.Range("D25").Select()
